{"js": "// Replace the three-digit-division problems in the worksheet table with a\n// new set of problems, preserving document order (including the one\n// duplicated prompt \"619\u00f74=\", whose two occurrences must map to two\n// different replacements).\nconst replacements = [\n  [\"373\u00f75=\", \"625\u00f79=\"],\n  [\"995\u00f79=\", \"484\u00f76=\"],\n  [\"276\u00f74=\", \"303\u00f75=\"],\n  [\"176\u00f73=\", \"655\u00f78=\"],\n  [\"408\u00f72=\", \"265\u00f73=\"],\n  [\"116\u00f72=\", \"781\u00f75=\"],\n  [\"869\u00f77=\", \"176\u00f76=\"],\n  [\"441\u00f72=\", \"375\u00f77=\"],\n  [\"631\u00f75=\", \"714\u00f72=\"],\n  [\"619\u00f74=\", \"612\u00f77=\"],\n  [\"826\u00f77=\", \"702\u00f73=\"],\n  [\"164\u00f78=\", \"904\u00f77=\"],\n  [\"453\u00f73=\", \"939\u00f72=\"],\n  [\"985\u00f76=\", \"354\u00f74=\"],\n  [\"619\u00f74=\", \"271\u00f76=\"],\n  [\"315\u00f77=\", \"632\u00f78=\"],\n  [\"163\u00f74=\", \"846\u00f77=\"],\n  [\"768\u00f78=\", \"328\u00f73=\"],\n  [\"748\u00f77=\", \"520\u00f78=\"],\n  [\"124\u00f77=\", \"880\u00f77=\"],\n  [\"146\u00f78=\", \"915\u00f79=\"],\n  [\"720\u00f77=\", \"380\u00f75=\"],\n  [\"953\u00f74=\", \"394\u00f72=\"],\n  [\"666\u00f79=\", \"426\u00f77=\"],\n  [\"567\u00f72=\", \"142\u00f78=\"],\n];\n\n// Group the (ordered) replacement list by the \"old\" text so that, for any\n// text value that appears more than once in the document, each occurrence\n// (in document order) is mapped to the correct replacement in the order the\n// pairs were declared above.\nconst byOld = new Map();\nfor (const [oldText, newText] of replacements) {\n  if (!byOld.has(oldText)) byOld.set(oldText, []);\n  byOld.get(oldText).push(newText);\n}\n\nconst body = context.document.body;\nconst searchResultsByOld = new Map();\nfor (const oldText of byOld.keys()) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"text\");\n  searchResultsByOld.set(oldText, found);\n}\nawait context.sync();\n\nfor (const [oldText, newTexts] of byOld) {\n  const found = searchResultsByOld.get(oldText);\n  for (let i = 0; i < found.items.length && i < newTexts.length; i++) {\n    found.items[i].insertText(newTexts[i], \"Replace\");\n  }\n}\nawait context.sync();\n", "ps1": "# Replace the three-digit-division problems in the worksheet table with a\n# new set of problems. Cells are addressed by (row, column) so the one\n# duplicated prompt (\"619\u00f74=\", appearing in row 5 col 5 and row 9 col 5)\n# is handled unambiguously - each occurrence gets its own replacement.\n$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n$replacements = @(\n    @{Row=1;  Col=1; Old=\"373\u00f75=\"; New=\"625\u00f79=\"},\n    @{Row=1;  Col=2; Old=\"995\u00f79=\"; New=\"484\u00f76=\"},\n    @{Row=1;  Col=3; Old=\"276\u00f74=\"; New=\"303\u00f75=\"},\n    @{Row=1;  Col=4; Old=\"176\u00f73=\"; New=\"655\u00f78=\"},\n    @{Row=1;  Col=5; Old=\"408\u00f72=\"; New=\"265\u00f73=\"},\n    @{Row=5;  Col=1; Old=\"116\u00f72=\"; New=\"781\u00f75=\"},\n    @{Row=5;  Col=2; Old=\"869\u00f77=\"; New=\"176\u00f76=\"},\n    @{Row=5;  Col=3; Old=\"441\u00f72=\"; New=\"375\u00f77=\"},\n    @{Row=5;  Col=4; Old=\"631\u00f75=\"; New=\"714\u00f72=\"},\n    @{Row=5;  Col=5; Old=\"619\u00f74=\"; New=\"612\u00f77=\"},\n    @{Row=9;  Col=1; Old=\"826\u00f77=\"; New=\"702\u00f73=\"},\n    @{Row=9;  Col=2; Old=\"164\u00f78=\"; New=\"904\u00f77=\"},\n    @{Row=9;  Col=3; Old=\"453\u00f73=\"; New=\"939\u00f72=\"},\n    @{Row=9;  Col=4; Old=\"985\u00f76=\"; New=\"354\u00f74=\"},\n    @{Row=9;  Col=5; Old=\"619\u00f74=\"; New=\"271\u00f76=\"},\n    @{Row=13; Col=1; Old=\"315\u00f77=\"; New=\"632\u00f78=\"},\n    @{Row=13; Col=2; Old=\"163\u00f74=\"; New=\"846\u00f77=\"},\n    @{Row=13; Col=3; Old=\"768\u00f78=\"; New=\"328\u00f73=\"},\n    @{Row=13; Col=4; Old=\"748\u00f77=\"; New=\"520\u00f78=\"},\n    @{Row=13; Col=5; Old=\"124\u00f77=\"; New=\"880\u00f77=\"},\n    @{Row=17; Col=1; Old=\"146\u00f78=\"; New=\"915\u00f79=\"},\n    @{Row=17; Col=2; Old=\"720\u00f77=\"; New=\"380\u00f75=\"},\n    @{Row=17; Col=3; Old=\"953\u00f74=\"; New=\"394\u00f72=\"},\n    @{Row=17; Col=4; Old=\"666\u00f79=\"; New=\"426\u00f77=\"},\n    @{Row=17; Col=5; Old=\"567\u00f72=\"; New=\"142\u00f78=\"}\n)\n\nforeach ($item in $replacements) {\n    $cell = $t.Cell($item.Row, $item.Col)\n    # Setting Range.Text replaces the cell's content in place while keeping\n    # the existing run/paragraph formatting (font, size, etc.) untouched.\n    $cell.Range.Text = $item.New\n}\n"}
